$d = $word.ActiveDocument

# The document has a single section whose default footer (footer1.xml) and
# first-page footer (footer2.xml) both carry the Pearson logo inline picture,
# while the first-page header (header2.xml) carries the BTEC logo inline
# picture. Word's InlineShape object has no settable "Name" — only the
# (floating) Shape object does — so round-trip each inline picture through
# ConvertToShape()/ConvertToInlineShape() to rename it while preserving its
# inline layout.
function Rename-InlinePicture($range, $newName) {
    $inlineShape = $range.InlineShapes.Item(1)
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

$sec = $d.Sections.Item(1)

# Footer (default) -> Pearson logo: image1.png -> image2.png
Rename-InlinePicture $sec.Footers.Item(1).Range "image2.png"

# Footer (first page) -> Pearson logo: image1.png -> image2.png
Rename-InlinePicture $sec.Footers.Item(2).Range "image2.png"

# Header (first page) -> BTEC logo: image2.jpg -> image1.jpg
Rename-InlinePicture $sec.Headers.Item(2).Range "image1.jpg"
